$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rates text in cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Text

$oldBlock = [char]0x2705 + " 1000 Bs = 1.81 = 6597.13 pesos`n" + [char]0x2705 + " 6597.13 pesos = 1.81 = 956.87 Bs"
$newBlock = [char]0x2705 + " 1000 Bs = 1.82 = 6635.61 pesos`n" + [char]0x2705 + " 6635.61 pesos = 1.81 = 955.72 Bs"

$cellA1.Value = $text.Replace($oldBlock, $newBlock)

# --- Sheet "tasas": update the rate values in N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 549
$ws2.Range("O10").Value = 3642.95
$ws2.Range("N12").Value = 3656.2
$ws2.Range("O12").Value = 526.6
